$wb = $excel.ActiveWorkbook

$newStamp = "2026-02-15 22:07"

# --- 1) PODSUMOWANIE: refresh the "last checked" timestamp for every profile row ---
$summary = $wb.Worksheets.Item("PODSUMOWANIE")
$summary.Range("B2").Value = $newStamp
$summary.Range("B3").Value = $newStamp
$summary.Range("B4").Value = $newStamp
$summary.Range("B5").Value = $newStamp
$summary.Range("B6").Value = $newStamp

# --- 2) Per-profile detail sheets: append a fresh monitoring-run row - row 13.
# B value per sheet = total listing count carried over unchanged from the
# previous run - no new/removed listings this pass; I value = the "details"
# note carried over for sheets that track individual ad codes.
$details = @(
    @{ Name = "wszystkie-lublin"; Count = 432; Notes = $null },
    @{ Name = "artymiuk";         Count = 0;   Notes = $null },
    @{ Name = "poqui";            Count = 5;   Notes = "1951OR|17NeTz|17vbYq|18KAEc|183ger" },
    @{ Name = "stylowepokoje";    Count = 2;   Notes = "195dLc|16ZeYm" },
    @{ Name = "villahome";        Count = 0;   Notes = $null }
)

foreach ($d in $details) {
    $ws = $wb.Worksheets.Item($d.Name)

    # Row 11 carries the "odd" banded-row styling that row 13 should continue;
    # row 12, the current last row, carries the alternate "even" styling.
    $ws.Range("A11:H11").Copy($ws.Range("A13:H13"))
    # Bring over column I's cell - present, possibly empty, on every data row -
    # from row 12 so the hidden notes column keeps its shape on the new row.
    $ws.Range("I12").Copy($ws.Range("I13"))
    $ws.Rows.Item(13).RowHeight = 18

    $ws.Range("A13").Value = $newStamp
    $ws.Range("B13").Value = $d.Count
    $ws.Range("C13").Value = 0
    $ws.Range("D13").Value = 0
    $ws.Range("E13").Value = 0
    $ws.Range("F13").Value = "—"
    $ws.Range("G13").Value = "—"
    $ws.Range("H13").Value = "OK"
    if ($d.Notes) {
        $ws.Range("I13").Value = $d.Notes
    }
    # else: I13 already holds the empty string carried over from I12's copy.

    # The "new last row" marker cell in column I only ever lived on the
    # previous last row - row 12 - when that column had no real content for
    # this sheet; now that row 13 is the last row, drop it from row 12.
    if (-not $d.Notes) {
        $ws.Range("I12").ClearContents()
    }
}
